# The commit moves the embedded scatter chart ("P-R_meanDetectionScore_9SP")
# to a new position on the sheet while keeping its overall size unchanged:
#   from: col 4 / +142875 EMU, row 5 / +123825 EMU   ->  col 4 / +247650 EMU, row 4 / +38100 EMU
#   to:   col 11 / +447675 EMU, row 20 / +9525 EMU    ->  col 11 / +552450 EMU, row 18 / +114300 EMU
#
# With the sheet's default column width (58.4375pt) and row height (15pt),
# those anchor cells correspond to absolute Left/Top of 245pt/84.75pt (before)
# and 253.25pt/63pt (after) - i.e. the chart was dragged up and to the right
# by the same amount on both axes, with Width (433.0625pt) and Height (216pt)
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$chartObj = $ws.ChartObjects(1)

# Keep size identical, only reposition (matches the anchor delta in the diff).
$chartObj.Left = 253.25
$chartObj.Top = 63
$chartObj.Width = 433.0625
$chartObj.Height = 216
